# Auto-generated Excel COM-interop script updating leve-profit market data cells
# across the ALC/ARM/BSM/CRP/CUL/GSM/LTW/WVR sheets (scheduled market-data refresh).

$wb = $excel.ActiveWorkbook

# --- ALC sheet ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H18").Value = 397.5
$ws.Range("I18").Value = 397.5
$ws.Range("K18").Value = 397.5
$ws.Range("M18").Value = -113.5
$ws.Range("H19").Value = 959.2222
$ws.Range("I19").Value = 1210
$ws.Range("J19").Value = 758.6
$ws.Range("K19").Value = 1210
$ws.Range("L19").Value = 758.6
$ws.Range("M19").Value = -1035
$ws.Range("N19").Value = -1108.6
$ws.Range("H28").Value = 51019.05
$ws.Range("I28").Value = 67190.734
$ws.Range("J28").Value = 2504
$ws.Range("K28").Value = 67190.734
$ws.Range("L28").Value = 2504
$ws.Range("M28").Value = -66705.734
$ws.Range("N28").Value = -3474
$ws.Range("H57").Value = 61780
$ws.Range("J57").Value = 61780
$ws.Range("L57").Value = 185340
$ws.Range("N57").Value = -186338
$ws.Range("H76").Value = 18250.5
$ws.Range("I76").Value = 17999.334
$ws.Range("J76").Value = 19004
$ws.Range("K76").Value = 17999.334
$ws.Range("L76").Value = 19004
$ws.Range("M76").Value = -17684.334
$ws.Range("N76").Value = -19634
$ws.Range("H79").Value = 18250.5
$ws.Range("I79").Value = 17999.334
$ws.Range("J79").Value = 19004
$ws.Range("K79").Value = 17999.334
$ws.Range("L79").Value = 19004
$ws.Range("M79").Value = -16907.334
$ws.Range("N79").Value = -21188
$ws.Range("H129").Value = 2465.25
$ws.Range("J129").Value = 3060.25
$ws.Range("L129").Value = 9180.75
$ws.Range("N129").Value = -19180.75
$ws.Range("H132").Value = 1889.7858
$ws.Range("I132").Value = 1650.0769
$ws.Range("K132").Value = 4950.2307
$ws.Range("M132").Value = -2420.2307

# --- ARM sheet ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 3057.1692
$ws.Range("I32").Value = 2651.5938
$ws.Range("K32").Value = 2651.5938
$ws.Range("M32").Value = -2364.5938
$ws.Range("H45").Value = 90914270
$ws.Range("I45").Value = 142858270
$ws.Range("J45").Value = 12253
$ws.Range("K45").Value = 142858270
$ws.Range("L45").Value = 12253
$ws.Range("M45").Value = -142857893
$ws.Range("N45").Value = -13007
$ws.Range("H61").Value = 10043.556
$ws.Range("I61").Value = 8284.714
$ws.Range("K61").Value = 8284.714
$ws.Range("M61").Value = -8072.714
$ws.Range("H132").Value = 3249.8
$ws.Range("I132").Value = 2600.2778
$ws.Range("K132").Value = 7800.8334
$ws.Range("M132").Value = -5270.8334
$ws.Range("H136").Value = 10043.556
$ws.Range("I136").Value = 8284.714
$ws.Range("K136").Value = 24854.142
$ws.Range("M136").Value = -22304.142

# --- BSM sheet ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H22").Value = 433.75
$ws.Range("I22").Value = 245
$ws.Range("J22").Value = 1000
$ws.Range("K22").Value = 245
$ws.Range("L22").Value = 1000
$ws.Range("M22").Value = -72
$ws.Range("N22").Value = -1346
$ws.Range("H80").Value = 1217
$ws.Range("J80").Value = 941.1429000000001
$ws.Range("L80").Value = 941.1429000000001
$ws.Range("N80").Value = -2937.1429
$ws.Range("H83").Value = 1217
$ws.Range("J83").Value = 941.1429000000001
$ws.Range("L83").Value = 4705.7145
$ws.Range("N83").Value = -14689.7145
$ws.Range("H134").Value = 1859.7142
$ws.Range("I134").Value = 1187.9474
$ws.Range("K134").Value = 3563.8422
$ws.Range("M134").Value = -1028.8422

# --- CRP sheet ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H5").Value = 1563.1818
$ws.Range("I5").Value = 598.5
$ws.Range("J5").Value = 1777.5555
$ws.Range("K5").Value = 598.5
$ws.Range("L5").Value = 1777.5555
$ws.Range("M5").Value = -486.5
$ws.Range("N5").Value = -2001.5555
$ws.Range("H8").Value = 2171.4285
$ws.Range("I8").Value = 200
$ws.Range("J8").Value = 2500
$ws.Range("K8").Value = 200
$ws.Range("L8").Value = 2500
$ws.Range("M8").Value = -60
$ws.Range("N8").Value = -2780
$ws.Range("H10").Value = 1592.8334
$ws.Range("I10").Value = 103.5
$ws.Range("J10").Value = 2337.5
$ws.Range("K10").Value = 103.5
$ws.Range("L10").Value = 2337.5
$ws.Range("M10").Value = 35.5
$ws.Range("N10").Value = -2615.5
$ws.Range("H11").Value = 462.33334
$ws.Range("I11").Value = 541.6667
$ws.Range("J11").Value = 383
$ws.Range("K11").Value = 541.6667
$ws.Range("L11").Value = 383
$ws.Range("M11").Value = -401.6667
$ws.Range("N11").Value = -663
$ws.Range("H14").Value = 2549.5
$ws.Range("I14").Value = 100
$ws.Range("J14").Value = 4999
$ws.Range("K14").Value = 100
$ws.Range("L14").Value = 4999
$ws.Range("M14").Value = 70
$ws.Range("N14").Value = -5339
$ws.Range("H15").Value = 4628.3706
$ws.Range("I15").Value = 3209.4211
$ws.Range("J15").Value = 7998.375
$ws.Range("K15").Value = 3209.4211
$ws.Range("L15").Value = 7998.375
$ws.Range("M15").Value = -3039.4211
$ws.Range("N15").Value = -8338.375
$ws.Range("H22").Value = 1217.3
$ws.Range("I22").Value = 399.42856
$ws.Range("J22").Value = 3125.6667
$ws.Range("K22").Value = 399.42856
$ws.Range("L22").Value = 3125.6667
$ws.Range("M22").Value = -49.42856
$ws.Range("N22").Value = -3825.6667
$ws.Range("H31").Value = 42642.57
$ws.Range("I31").Value = 4711.2856
$ws.Range("J31").Value = 80573.86
$ws.Range("K31").Value = 4711.2856
$ws.Range("L31").Value = 80573.86
$ws.Range("M31").Value = -4416.2856
$ws.Range("N31").Value = -81163.86
$ws.Range("H34").Value = 42642.57
$ws.Range("I34").Value = 4711.2856
$ws.Range("J34").Value = 80573.86
$ws.Range("K34").Value = 4711.2856
$ws.Range("L34").Value = 80573.86
$ws.Range("M34").Value = -4509.2856
$ws.Range("N34").Value = -80977.86
$ws.Range("H105").Value = 6908.091
$ws.Range("I105").Value = 4593.625
$ws.Range("J105").Value = 13080
$ws.Range("K105").Value = 4593.625
$ws.Range("L105").Value = 13080
$ws.Range("M105").Value = -2846.625
$ws.Range("N105").Value = -16574
$ws.Range("H107").Value = 908.52
$ws.Range("I107").Value = 801.2222
$ws.Range("J107").Value = 1184.4286
$ws.Range("K107").Value = 801.2222
$ws.Range("L107").Value = 1184.4286
$ws.Range("M107").Value = 1118.7778
$ws.Range("N107").Value = -5024.4286
$ws.Range("H132").Value = 10309.692
$ws.Range("I132").Value = 11445.777
$ws.Range("J132").Value = 7753.5
$ws.Range("K132").Value = 34337.331
$ws.Range("L132").Value = 23260.5
$ws.Range("M132").Value = -31807.331
$ws.Range("N132").Value = -28320.5
$ws.Range("H141").Value = 122919.57
$ws.Range("J141").Value = 122919.57
$ws.Range("L141").Value = 122919.57
$ws.Range("N141").Value = -133279.57

# --- CUL sheet ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H32").Value = 20906
$ws.Range("I32").Value = 490
$ws.Range("J32").Value = 34516.668
$ws.Range("K32").Value = 1470
$ws.Range("L32").Value = 103550.004
$ws.Range("M32").Value = -1187
$ws.Range("N32").Value = -104116.004
$ws.Range("H59").Value = 1400
$ws.Range("J59").Value = 500
$ws.Range("L59").Value = 1500
$ws.Range("N59").Value = -2580
$ws.Range("H60").Value = 296.55554
$ws.Range("I60").Value = 296.75
$ws.Range("J60").Value = 295
$ws.Range("K60").Value = 890.25
$ws.Range("L60").Value = 885
$ws.Range("M60").Value = -639.25
$ws.Range("N60").Value = -1387
$ws.Range("H114").Value = 2000
$ws.Range("I114").Value = 500
$ws.Range("J114").Value = 5000
$ws.Range("K114").Value = 1500
$ws.Range("L114").Value = 15000
$ws.Range("M114").Value = 1754
$ws.Range("N114").Value = -21508
$ws.Range("H139").Value = 5538.222
$ws.Range("I139").Value = 2229.6667
$ws.Range("J139").Value = 12155.333
$ws.Range("K139").Value = 6689.000100000001
$ws.Range("L139").Value = 36465.999
$ws.Range("M139").Value = -1549.000100000001
$ws.Range("N139").Value = -46745.999

# --- GSM sheet ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H113").Value = 3560.8
$ws.Range("I113").Value = 2536.5
$ws.Range("K113").Value = 2536.5
$ws.Range("M113").Value = -366.5
$ws.Range("H122").Value = 4272.4346
$ws.Range("I122").Value = 3980.0715
$ws.Range("K122").Value = 11940.2145
$ws.Range("M122").Value = -9490.2145

# --- LTW sheet ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H16").Value = 2055.25
$ws.Range("I16").Value = 2040.5333
$ws.Range("K16").Value = 2040.5333
$ws.Range("M16").Value = -1870.5333
$ws.Range("H40").Value = 9622.684999999999
$ws.Range("I40").Value = 9844.357
$ws.Range("J40").Value = 9002
$ws.Range("K40").Value = 9844.357
$ws.Range("L40").Value = 9002
$ws.Range("M40").Value = -9708.357
$ws.Range("N40").Value = -9274
$ws.Range("H82").Value = 6356.7856
$ws.Range("I82").Value = 6712.5
$ws.Range("K82").Value = 6712.5
$ws.Range("M82").Value = -6351.5
$ws.Range("H85").Value = 6356.7856
$ws.Range("I85").Value = 6712.5
$ws.Range("K85").Value = 6712.5
$ws.Range("M85").Value = -5464.5

# --- WVR sheet ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H60").Value = 0
$ws.Range("J60").Value = 0
$ws.Range("L60").Value = 0
$ws.Range("N60").Value = ""
$ws.Range("H62").Value = 7940.3
$ws.Range("I62").Value = 7280.8
$ws.Range("J62").Value = 8599.799999999999
$ws.Range("K62").Value = 7280.8
$ws.Range("L62").Value = 8599.799999999999
$ws.Range("M62").Value = -6656.8
$ws.Range("N62").Value = -9847.799999999999
$ws.Range("H65").Value = 7940.3
$ws.Range("I65").Value = 7280.8
$ws.Range("J65").Value = 8599.799999999999
$ws.Range("K65").Value = 36404
$ws.Range("L65").Value = 42999
$ws.Range("M65").Value = -33284
$ws.Range("N65").Value = -49239
$ws.Range("H74").Value = 14248
$ws.Range("J74").Value = 14248
$ws.Range("L74").Value = 14248
$ws.Range("N74").Value = -16120
$ws.Range("H77").Value = 14248
$ws.Range("J77").Value = 14248
$ws.Range("L77").Value = 42744
$ws.Range("N77").Value = -52104
$ws.Range("H81").Value = 3693.9167
$ws.Range("I81").Value = 2477.8572
$ws.Range("K81").Value = 4955.7144
$ws.Range("M81").Value = -3894.7144
$ws.Range("H84").Value = 3693.9167
$ws.Range("I84").Value = 2477.8572
$ws.Range("K84").Value = 24778.572
$ws.Range("M84").Value = -19474.572
